$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (existing anomaly record) - desvio_padrao / media_vendas adjusted
$ws.Range("H5").Value = 1.06
$ws.Range("I5").Value = 0.24

# Add new row 7 with the new atypical sale record for 2025-06-14.
# New shared-string entries are appended in the order they are first
# assigned, so write the cells that introduce brand-new strings first
# and in the same order as the source export (date, cliente, id_venda),
# matching how the reference file's shared string table was built.
#
# A7 and D7 must stay plain text (not auto-converted to date/number),
# so force text format, assign, then clear the format override so no
# style index is left on the cell (matches the other data rows).
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-06-14"
$ws.Range("A7").ClearFormats()

$ws.Range("C7").Value = "SOCIEDADE MICHELIN DE PARTICIPACOES INDUST E COMERCIO LTDA"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "360346"
$ws.Range("D7").ClearFormats()

$ws.Range("B7").Value = 2
$ws.Range("E7").Value = 13588
$ws.Range("F7").Value = "CANETA STYLLUS ACTIVA AGOLD"
$ws.Range("G7").Value = -4
$ws.Range("H7").Value = 1.06
$ws.Range("I7").Value = 0.24
